$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.016.93"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -1.53%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.640.26"
$cell.ClearFormats()
$ws.Range("E3").Value = "  -1.69%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.015"
$cell.ClearFormats()
$ws.Range("E4").Value = "  +0.73%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "216.33"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5016"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.015"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +0.75%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2580"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +0.02%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06441"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -0.33%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "19.50"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -2.59%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07771"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +1.37%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.643.51"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -1.60%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.266"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -1.88%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.865.31"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -1.57%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.5466"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -1.71%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0₅7948"
$cell.ClearFormats()
$ws.Range("E16").Value = "  -1.27%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "63.73"
$cell.ClearFormats()
$ws.Range("E17").Value = "  -1.62%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "26.038.66"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -1.45%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "1.017"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +0.85%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "204.44"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -2.67%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.319"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -2.45%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.02"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -1.03%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.980"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  +0.82%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.970"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +13.91%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "141.55"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -2.56%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1153"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -1.28%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "15.82"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +0.10%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "6.814"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -2.87%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.05052"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -3.59%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.245"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -1.49%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.271"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -3.22%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.211"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -0.41%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.549"
$cell.ClearFormats()
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.354"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -1.01%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.8947"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -3.79%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.622"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -4.86%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.5655"
$cell.ClearFormats()
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.132.98"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -1.85%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01564"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("E41").Value = "  +0.48%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.016"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +0.76%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.646"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -0.12%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.8190"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -3.25%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "99.94"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.45%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.774.99"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("E47").Value = "  +2.34%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.4549"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +1.26%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.018"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +1.29%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "54.94"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -1.94%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.05035"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -1.54%  "
